$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Insert the large block of new to-do list paragraphs right before
# the "NEW KEYWORD IN CODE" paragraph (i.e. right after the paragraph that
# ends with "COMMENTCOMMENTCOMMENT + Class diagram.").
# ---------------------------------------------------------------------------

$anchor = $d.Content
$anchor.Find.ClearFormatting()
$foundAnchor = $anchor.Find.Execute("COMMENTCOMMENTCOMMENT + Class diagram.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundAnchor) {
    throw "Could not find anchor paragraph 'COMMENTCOMMENTCOMMENT + Class diagram.'"
}

# Collapse to the end of the found text, then extend to the end of the
# paragraph (past the paragraph mark) so the new content lands in its own
# paragraph, immediately before the "NEW KEYWORD IN CODE" paragraph.
$insertRange = $anchor.Duplicate
$insertRange.Collapse(0)
$insertRange.MoveEndUntil("X", 1) | Out-Null
$insertRange.MoveEnd(5, 1) | Out-Null

$newParasXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="11"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t>Weapons and Attacking</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="11"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t>Plan of attack</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="2"/>
      <w:numId w:val="11"/>
    </w:numPr>
    <w:rPr>
      <w:strike/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:strike/>
    </w:rPr>
    <w:t xml:space="preserve">Add the following to </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:strike/>
    </w:rPr>
    <w:t>RoguelikeObject</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:strike/>
    </w:rPr>
    <w:t xml:space="preserve">: </w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="3"/>
      <w:numId w:val="11"/>
    </w:numPr>
    <w:rPr>
      <w:strike/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:strike/>
    </w:rPr>
    <w:t>Melee Damage</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="3"/>
      <w:numId w:val="11"/>
    </w:numPr>
    <w:rPr>
      <w:strike/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:strike/>
    </w:rPr>
    <w:t>Thrown Damage</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="3"/>
      <w:numId w:val="11"/>
    </w:numPr>
    <w:rPr>
      <w:strike/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:strike/>
    </w:rPr>
    <w:t>Attack Speed</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="2"/>
      <w:numId w:val="11"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t>Add the following to Weapon:</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="3"/>
      <w:numId w:val="11"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:lastRenderedPageBreak/>
    <w:t>Unique Context Menu</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="3"/>
      <w:numId w:val="11"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t>Unique Log Entries</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="2"/>
      <w:numId w:val="11"/>
    </w:numPr>
    <w:rPr>
      <w:strike/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:strike/>
    </w:rPr>
    <w:t>Add the following to Equipment:</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="3"/>
      <w:numId w:val="11"/>
    </w:numPr>
    <w:rPr>
      <w:strike/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:strike/>
    </w:rPr>
    <w:t>Unique &#8220;Equip&#8221; Context menu</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="2"/>
      <w:numId w:val="11"/>
    </w:numPr>
    <w:rPr>
      <w:strike/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:strike/>
    </w:rPr>
    <w:t xml:space="preserve">Units </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:strike/>
    </w:rPr>
    <w:t>have list of gibs which like 2 legs, 2 arms, a torso and a head</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="3"/>
      <w:numId w:val="11"/>
    </w:numPr>
    <w:rPr>
      <w:strike/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:strike/>
    </w:rPr>
    <w:t>Some gibs have the grasping bool</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="2"/>
      <w:numId w:val="11"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t>Clean up what was added above ^^^</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="2"/>
      <w:numId w:val="11"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t>Add a function for quickly adding a body part and a function for removing a bodypart.</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$insertRange.InsertXML($newParasXml)

# ---------------------------------------------------------------------------
# Step 2: Merge the split "...removed in the f" + bookmark + "uture." runs
# into a single run now that the bookmark has moved to the new paragraph
# above, while keeping the surrounding tab / [ / REMOVE? / ] / br / [NEEDS
# WORK] runs intact and in the same order.
# ---------------------------------------------------------------------------

$findRange = $d.Content
$findRange.Find.ClearFormatting()
$foundTarget = $findRange.Find.Execute("[REMOVE?]", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundTarget) {
    throw "Could not find '[REMOVE?]' marker"
}

$replaceRange = $findRange.Duplicate
$replaceRange.Collapse(1)
$replaceRange.MoveEndUntil([char]9, 1) | Out-Null

$replacementXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:r><w:t>[</w:t></w:r>
  <w:r><w:t>REMOVE?</w:t></w:r>
  <w:r><w:t>]</w:t></w:r>
  <w:r><w:t xml:space="preserve"> &#8211; add this as a comment to code that may need to be removed in the future.</w:t></w:r>
  <w:r><w:br/></w:r>
  <w:r><w:tab/><w:t>[NEEDS WORK] &#8211; This function will need to be updated in the future</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$replaceRange.InsertXML($replacementXml)

# ---------------------------------------------------------------------------
# Step 3: Remove the stray <w:lastRenderedPageBreak/> from the "Takes an
# item and adds as much of it as possible to an inventory." paragraph.
# ---------------------------------------------------------------------------

$takeRange = $d.Content
$takeRange.Find.ClearFormatting()
$foundTake = $takeRange.Find.Execute("Takes an item and adds as much of it as possible to an inventory.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundTake) {
    throw "Could not find 'Takes an item and adds...' paragraph"
}

$takeParaRange = $takeRange.Paragraphs(1).Range.Duplicate

$takeReplacementXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="21"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t>Takes an item and adds as much of it as possible to an inventory.</w:t>
  </w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$takeParaRange.InsertXML($takeReplacementXml)

Write-Output "All edits applied."
